$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The transition from Max -> Model Viewer is now seamless, so the two
# "Scene Exporter" backlog items that used to track that work
# ("Make exporter a GUP...") and ("Make the path from Max->Model Viewer
# seamless") are no longer needed. Capture the two existing cell
# comments first (they sit a couple of rows further down the sheet, on
# the "Compiled shaders" / "Error handling strategy..." rows) so they
# can be re-anchored once those rows shift up.
$cmShaders = $ws.Range("B12").Comment
$textShaders = $cmShaders.Text()
$cmErrorHandling = $ws.Range("B16").Comment
$textErrorHandling = $cmErrorHandling.Text()
$cmShaders.Delete()
$cmErrorHandling.Delete()

# Remove the two obsolete rows entirely; everything below shifts up by two.
$ws.Range("A2:C3").EntireRow.Delete()

# Re-create the comments on their new (shifted up by two) cells.
$ws.Range("B10").AddComment($textShaders) | Out-Null
$ws.Range("B14").AddComment($textErrorHandling) | Out-Null

# Leave the selection the way it was after removing the two rows.
$ws.Range("A2:XFD3").Select()
